$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 93, shifting the existing rows 93-119 down to 95-121.
$ws.Rows("93:94").Insert()

# Row 93 (new weekly record)
$ws.Range("A93").Value = 10
$ws.Range("B93").Value = "Vega Modelo de Temuco"
$ws.Range("C93").Value = "La Araucanía"
$ws.Range("D93").Value = 44466
$ws.Range("E93").Value = 9
$ws.Range("F93").Value = 100112013
$ws.Range("G93").Value = "Alcachofa"
$ws.Range("H93").Value = "Española"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 180
$ws.Range("K93").Value = 12000
$ws.Range("L93").Value = 13000
$ws.Range("M93").Value = 12556
$ws.Range("N93").Value = "`$/caja 30 unidades"
$ws.Range("O93").Value = "Región Metropolitana"
$ws.Range("P93").Value = 419
$ws.Range("Q93").Value = 30
$ws.Range("R93").Value = "Hortaliza"

# Row 94 (new weekly record)
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 44466
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 100112013
$ws.Range("G94").Value = "Alcachofa"
$ws.Range("H94").Value = "Madrigal"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 170
$ws.Range("K94").Value = 12000
$ws.Range("L94").Value = 13000
$ws.Range("M94").Value = 12471
$ws.Range("N94").Value = "`$/caja 40 unidades"
$ws.Range("O94").Value = "Región Metropolitana"
$ws.Range("P94").Value = 312
$ws.Range("Q94").Value = 40
$ws.Range("R94").Value = "Hortaliza"

Write-Output "done"
